$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "String Column"
$ws.Range("B1").Value = "Integer Column"
$ws.Range("C1").Value = "Date Column"
$ws.Range("D1").Value = "Optional Column"
$ws.Range("E1").Value = "Optional Value Column"

$ws.Range("C2:C3").NumberFormat = "@"

$ws.Range("A2").Value = "this is a text"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "2016.12.09"
$ws.Range("D2").Value = "has value"

$ws.Range("A3").Value = "this is another text"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "2016.12.31"
$ws.Range("D3").Value = "has another value"
$ws.Range("E3").Value = "has value"

$ws.Columns.Item(1).ColumnWidth = 17.85546875
$ws.Columns.Item(2).ColumnWidth = 15.7109375
$ws.Columns.Item(3).ColumnWidth = 12.42578125
$ws.Columns.Item(4).ColumnWidth = 18.140625
$ws.Columns.Item(5).ColumnWidth = 22

$ws.Range("D3").Select()
